$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.679.17"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "1.643.82"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.529"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.36"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").Value = "1.876.96"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "1.644.54"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.68"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "27.663.29"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("E23").Value = "  +7.28%  "

$ws.Range("E24").Value = "  -2.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").Value = "1.444.64"
$ws.Range("E33").Value = "  +2.74%  "

$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.889"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.63%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.63"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.16%  "

$ws.Range("B45").Value = "mCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.54%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.786.68"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.81%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.90%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.06%  "
